$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a new list item "budget" right before the existing "category" list
#    item that follows the "BudgetCategory" bullet (ilvl 1, numId 4).
# ---------------------------------------------------------------------------
$targetText = "category"
$found = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    if ($txt -eq $targetText) {
        $prevTxt = $d.Paragraphs($i - 1).Range.Text
        $prevTxt = $prevTxt.TrimEnd([char]13, [char]7)
        if ($prevTxt -eq "BudgetCategory") {
            $found = $i
            break
        }
    }
}

$categoryPara = $d.Paragraphs($found)
$r = $categoryPara.Range
$r.Collapse(1)
$r.InsertParagraphBefore()
$newPara = $d.Paragraphs($found)
$newPara.Range.InsertBefore("budget")

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: it currently sits on the
#    "Ruleset" run and needs to move to the "Period end" run that precedes
#    it (paragination shifted because of the text added above).
# ---------------------------------------------------------------------------
function Find-ParaIndexAfter($doc, $text, $prevText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            if ($i -gt 1) {
                $p = $doc.Paragraphs($i - 1).Range.Text
                $p = $p.TrimEnd([char]13, [char]7)
                if ($p -eq $prevText) {
                    return $i
                }
            }
        }
    }
    return -1
}

$periodEndIdx = Find-ParaIndexAfter $d "Period end" "Period start"
$rulesetIdx = Find-ParaIndexAfter $d "Ruleset" "Period end"

# Remove the break from the "Ruleset" run by retyping its text - this keeps
# the paragraph's own pPr/rPr and w14:paraId/rsid attributes untouched while
# dropping the stray <w:lastRenderedPageBreak/> run child.
$rulesetPara = $d.Paragraphs($rulesetIdx)
$rulesetPara.Range.Text = "Ruleset"

# Re-fetch (structural edit may have reseated the range/collection).
$periodEndIdx = Find-ParaIndexAfter $d "Period end" "Period start"
$periodEndPara = $d.Paragraphs($periodEndIdx)
$peRange = $periodEndPara.Range
$peRange.Collapse(1)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="795048F3" w14:textId="63C6C59F" w:rsidR="00334FCE" w:rsidRDefault="00334FCE" w:rsidP="00896E4E"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="200" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:lang w:val="en"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:lang w:val="en"/></w:rPr><w:lastRenderedPageBreak/><w:t>Period end</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$peRange.InsertXML($xml)

Write-Host "Done."
